$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Płetwal błękitny" -> "Płetwal\nblekitny" (newline + wrap text) for the blue-whale row (B12)
$cell = $ws.Cells.Item(12, 2)
$cell.Value2 = "Płetwal" + [char]10 + "blekitny"
$cell.WrapText = $true

$ws.Range("B13").Select() | Out-Null
